$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "timespan" column (F) values from "Week" to "Month" for all data rows (2-18)
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Value2 -eq "Week") {
        $cell.Value2 = "Month"
    }
}

# Update the selected cell in the sheet view from J11 to J10
$ws.Range("J10").Select()
